$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 32.93949833333333
$ws.Range("H2").Value = 98.818495
$ws.Range("I2").Value = 0.02571831923682078
$ws.Range("J2").Value = 0.02571831923682077
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1451143333333333
$ws.Range("N2").Value = 0.435343
$ws.Range("O2").Value = 0.140827770705994
$ws.Range("P2").Value = 0.1408277707059941
$ws.Range("Q2").Value = 4.779993340976111
$ws.Range("R2").Value = 43.019940068785
$ws.Range("S2").Value = 0.003621853564426552
$ws.Range("T2").Value = 0.003621853564426552
$ws.Range("G3").Value = 32.93949833333333
$ws.Range("H3").Value = 98.818495
$ws.Range("I3").Value = 0.02571831923682078
$ws.Range("J3").Value = 0.02571831923682077
$ws.Range("O3").Value = 0.0422072807203407
$ws.Range("P3").Value = 0.0422072807203407
$ws.Range("Q3").Value = 1.432604661513333
$ws.Range("R3").Value = 12.89344195362
$ws.Range("S3").Value = 0.001085500319683833
$ws.Range("T3").Value = 0.001085500319683833
$ws.Range("G4").Value = 32.93949833333333
$ws.Range("H4").Value = 98.818495
$ws.Range("I4").Value = 0.02571831923682078
$ws.Range("J4").Value = 0.02571831923682077
$ws.Range("M4").Value = 0.841832
$ws.Range("N4").Value = 2.525496
$ws.Range("O4").Value = 0.8169649485736653
$ws.Range("P4").Value = 0.8169649485736653
$ws.Range("Q4").Value = 27.72952376094667
$ws.Range("R4").Value = 249.56571384852
$ws.Range("S4").Value = 0.0210109653527104
$ws.Range("T4").Value = 0.02101096535271039
$ws.Range("I5").Value = 0.4140443484779395
$ws.Range("J5").Value = 0.4140443484779395
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1451143333333333
$ws.Range("N5").Value = 0.435343
$ws.Range("O5").Value = 0.140827770705994
$ws.Range("P5").Value = 0.1408277707059941
$ws.Range("Q5").Value = 76.95406571358809
$ws.Range("R5").Value = 692.5865914222929
$ws.Range("S5").Value = 0.05830894256956395
$ws.Range("T5").Value = 0.05830894256956396
$ws.Range("I6").Value = 0.4140443484779395
$ws.Range("J6").Value = 0.4140443484779395
$ws.Range("O6").Value = 0.0422072807203407
$ws.Range("P6").Value = 0.0422072807203407
$ws.Range("S6").Value = 0.01747568604687896
$ws.Range("T6").Value = 0.01747568604687896
$ws.Range("I7").Value = 0.4140443484779395
$ws.Range("J7").Value = 0.4140443484779395
$ws.Range("M7").Value = 0.841832
$ws.Range("N7").Value = 2.525496
$ws.Range("O7").Value = 0.8169649485736653
$ws.Range("P7").Value = 0.8169649485736653
$ws.Range("Q7").Value = 446.4231310562107
$ws.Range("R7").Value = 4017.808179505896
$ws.Range("S7").Value = 0.3382597198614966
$ws.Range("T7").Value = 0.3382597198614966
$ws.Range("G8").Value = 422.1807963333333
$ws.Range("H8").Value = 1266.542389
$ws.Range("I8").Value = 0.3296279860087694
$ws.Range("J8").Value = 0.3296279860087693
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1451143333333333
$ws.Range("N8").Value = 0.435343
$ws.Range("O8").Value = 0.140827770705994
$ws.Range("P8").Value = 0.1408277707059941
$ws.Range("Q8").Value = 61.26448480604743
$ws.Range("R8").Value = 551.3803632544269
$ws.Range("S8").Value = 0.04642077443192158
$ws.Range("T8").Value = 0.04642077443192158
$ws.Range("G9").Value = 422.1807963333333
$ws.Range("H9").Value = 1266.542389
$ws.Range("I9").Value = 0.3296279860087694
$ws.Range("J9").Value = 0.3296279860087693
$ws.Range("O9").Value = 0.0422072807203407
$ws.Range("P9").Value = 0.0422072807203407
$ws.Range("Q9").Value = 18.36148719412934
$ws.Range("R9").Value = 165.253384747164
$ws.Range("S9").Value = 0.01391270093875266
$ws.Range("T9").Value = 0.01391270093875266
$ws.Range("G10").Value = 422.1807963333333
$ws.Range("H10").Value = 1266.542389
$ws.Range("I10").Value = 0.3296279860087694
$ws.Range("J10").Value = 0.3296279860087693
$ws.Range("M10").Value = 0.841832
$ws.Range("N10").Value = 2.525496
$ws.Range("O10").Value = 0.8169649485736653
$ws.Range("P10").Value = 0.8169649485736653
$ws.Range("Q10").Value = 355.4053041388827
$ws.Range("R10").Value = 3198.647737249944
$ws.Range("S10").Value = 0.2692945106380951
$ws.Range("T10").Value = 0.2692945106380951
$ws.Range("G11").Value = 16.509264
$ws.Range("H11").Value = 49.527792
$ws.Range("I11").Value = 0.01289001179132366
$ws.Range("J11").Value = 0.01289001179132366
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1451143333333333
$ws.Range("N11").Value = 0.435343
$ws.Range("O11").Value = 0.140827770705994
$ws.Range("P11").Value = 0.1408277707059941
$ws.Range("Q11").Value = 2.395730839184
$ws.Range("R11").Value = 21.561577552656
$ws.Range("S11").Value = 0.001815271624946088
$ws.Range("T11").Value = 0.001815271624946088
$ws.Range("G12").Value = 16.509264
$ws.Range("H12").Value = 49.527792
$ws.Range("I12").Value = 0.01289001179132366
$ws.Range("J12").Value = 0.01289001179132366
$ws.Range("O12").Value = 0.0422072807203407
$ws.Range("P12").Value = 0.0422072807203407
$ws.Range("Q12").Value = 0.7180209098879999
$ws.Range("R12").Value = 6.462188188992
$ws.Range("S12").Value = 0.0005440523461648993
$ws.Range("T12").Value = 0.0005440523461648994
$ws.Range("G13").Value = 16.509264
$ws.Range("H13").Value = 49.527792
$ws.Range("I13").Value = 0.01289001179132366
$ws.Range("J13").Value = 0.01289001179132366
$ws.Range("M13").Value = 0.841832
$ws.Range("N13").Value = 2.525496
$ws.Range("O13").Value = 0.8169649485736653
$ws.Range("P13").Value = 0.8169649485736653
$ws.Range("Q13").Value = 13.898026731648
$ws.Range("R13").Value = 125.082240584832
$ws.Range("S13").Value = 0.01053068782021267
$ws.Range("T13").Value = 0.01053068782021267
$ws.Range("G14").Value = 44.62094166666666
$ws.Range("H14").Value = 133.862825
$ws.Range("I14").Value = 0.03483889192294087
$ws.Range("J14").Value = 0.03483889192294087
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1451143333333333
$ws.Range("N14").Value = 0.435343
$ws.Range("O14").Value = 0.140827770705994
$ws.Range("P14").Value = 0.1408277707059941
$ws.Range("Q14").Value = 6.475138202663887
$ws.Range("R14").Value = 58.27624382397499
$ws.Range("S14").Value = 0.004906283483374825
$ws.Range("T14").Value = 0.004906283483374825
$ws.Range("G15").Value = 44.62094166666666
$ws.Range("H15").Value = 133.862825
$ws.Range("I15").Value = 0.03483889192294087
$ws.Range("J15").Value = 0.03483889192294087
$ws.Range("O15").Value = 0.0422072807203407
$ws.Range("P15").Value = 0.0422072807203407
$ws.Range("Q15").Value = 1.940653994966667
$ws.Range("R15").Value = 17.4658859547
$ws.Range("S15").Value = 0.001470454891377176
$ws.Range("T15").Value = 0.001470454891377176
$ws.Range("G16").Value = 44.62094166666666
$ws.Range("H16").Value = 133.862825
$ws.Range("I16").Value = 0.03483889192294087
$ws.Range("J16").Value = 0.03483889192294087
$ws.Range("M16").Value = 0.841832
$ws.Range("N16").Value = 2.525496
$ws.Range("O16").Value = 0.8169649485736653
$ws.Range("P16").Value = 0.8169649485736653
$ws.Range("Q16").Value = 37.56333656513333
$ws.Range("R16").Value = 338.0700290862
$ws.Range("S16").Value = 0.02846215354818887
$ws.Range("T16").Value = 0.02846215354818887
$ws.Range("G17").Value = 234.229538
$ws.Range("H17").Value = 702.6886139999999
$ws.Range("I17").Value = 0.1828804425622059
$ws.Range("J17").Value = 0.1828804425622059
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1451143333333333
$ws.Range("N17").Value = 0.435343
$ws.Range("O17").Value = 0.140827770705994
$ws.Range("P17").Value = 0.1408277707059941
$ws.Range("Q17").Value = 33.99006325384466
$ws.Range("R17").Value = 305.9105692846019
$ws.Range("S17").Value = 0.02575464503176104
$ws.Range("T17").Value = 0.02575464503176104
$ws.Range("G18").Value = 234.229538
$ws.Range("H18").Value = 702.6886139999999
$ws.Range("I18").Value = 0.1828804425622059
$ws.Range("J18").Value = 0.1828804425622059
$ws.Range("O18").Value = 0.0422072807203407
$ws.Range("P18").Value = 0.0422072807203407
$ws.Range("Q18").Value = 10.187111066696
$ws.Range("R18").Value = 91.683999600264
$ws.Range("S18").Value = 0.007718886177483166
$ws.Range("T18").Value = 0.007718886177483167
$ws.Range("G19").Value = 234.229538
$ws.Range("H19").Value = 702.6886139999999
$ws.Range("I19").Value = 0.1828804425622059
$ws.Range("J19").Value = 0.1828804425622059
$ws.Range("M19").Value = 0.841832
$ws.Range("N19").Value = 2.525496
$ws.Range("O19").Value = 0.8169649485736653
$ws.Range("P19").Value = 0.8169649485736653
$ws.Range("Q19").Value = 197.181920433616
$ws.Range("R19").Value = 1774.637283902544
$ws.Range("S19").Value = 0.1494069113529617
$ws.Range("T19").Value = 0.1494069113529617

Write-Host "Applied 208 cell updates"